$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data below the existing content
$ws.Range("B24").Value = 0.5
$ws.Range("C24").Value = "storypoint pr time"

# Update the selected cell shown in the UI
$ws.Range("M5").Select()
